$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells that would otherwise be auto-parsed as numbers
$textCells = @("D5", "D6", "D7", "D8", "D10", "D11", "D12", "D16", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D29", "D31", "D32", "D33", "D35", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "54.787.26"
$ws.Range("E2").Value = "  +6.55%  "

$ws.Range("D3").Value = "2.431.83"
$ws.Range("E3").Value = "  +7.03%  "

$ws.Range("E4").Value = "  +0.82%  "

$ws.Range("D5").Value = "475.37"
$ws.Range("E5").Value = "  +10.67%  "

$ws.Range("D6").Value = "138.37"
$ws.Range("E6").Value = "  +16.70%  "

$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  +0.46%  "

$ws.Range("D8").Value = "0.500"
$ws.Range("E8").Value = "  +10.60%  "

$ws.Range("D9").Value = "2.450.27"
$ws.Range("E9").Value = "  +8.15%  "

$ws.Range("D10").Value = "0.0954"
$ws.Range("E10").Value = "  +9.58%  "

$ws.Range("D11").Value = "5.49"
$ws.Range("E11").Value = "  +9.58%  "

$ws.Range("D12").Value = "0.323"
$ws.Range("E12").Value = "  +8.91%  "

$ws.Range("E13").Value = "  +3.14%  "

$ws.Range("D14").Value = "2.845.84"
$ws.Range("E14").Value = "  +8.17%  "

$ws.Range("D15").Value = "54.914.94"
$ws.Range("E15").Value = "  +7.30%  "

$ws.Range("D16").Value = "20.27"
$ws.Range("E16").Value = "  +10.27%  "

$ws.Range("E17").Value = "  +14.78%  "

$ws.Range("D18").Value = "2.439.64"
$ws.Range("E18").Value = "  +10.30%  "

$ws.Range("D19").Value = "4.32"
$ws.Range("E19").Value = "  +8.99%  "

$ws.Range("D20").Value = "9.81"
$ws.Range("E20").Value = "  +14.76%  "

$ws.Range("D21").Value = "311.56"
$ws.Range("E21").Value = "  +6.61%  "

$ws.Range("D22").Value = "0.991"
$ws.Range("E22").Value = "  -2.12%  "

$ws.Range("D23").Value = "5.65"
$ws.Range("E23").Value = "  +11.70%  "

$ws.Range("D24").Value = "56.84"
$ws.Range("E24").Value = "  +9.02%  "

$ws.Range("D25").Value = "0.998"
$ws.Range("E25").Value = "  -0.49%  "

$ws.Range("D26").Value = "0.399"
$ws.Range("E26").Value = "  +9.55%  "

$ws.Range("E27").Value = "  +23.16%  "

$ws.Range("D28").Value = "2.533.79"
$ws.Range("E28").Value = "  +9.72%  "

$ws.Range("D29").Value = "7.25"
$ws.Range("E29").Value = "  +9.70%  "

$ws.Range("D30").Value = "0.0₃0762"
$ws.Range("E30").Value = "  +17.23%  "

$ws.Range("D31").Value = "0.998"
$ws.Range("E31").Value = "  +0.47%  "

$ws.Range("D32").Value = "147.60"
$ws.Range("E32").Value = "  +4.03%  "

$ws.Range("D33").Value = "18.00"
$ws.Range("E33").Value = "  +9.13%  "

$ws.Range("E34").Value = "  +11.64%  "

$ws.Range("D35").Value = "5.08"
$ws.Range("E35").Value = "  +9.75%  "

$ws.Range("E36").Value = "  +15.57%  "

$ws.Range("D37").Value = "3.55"
$ws.Range("E37").Value = "  +9.79%  "

$ws.Range("D38").Value = "0.828"
$ws.Range("E38").Value = "  +11.50%  "

$ws.Range("D39").Value = "33.66"
$ws.Range("E39").Value = "  +5.91%  "

$ws.Range("D40").Value = "0.993"
$ws.Range("E40").Value = "  +0.94%  "

$ws.Range("D41").Value = "3.41"
$ws.Range("E41").Value = "  +10.39%  "

$ws.Range("D42").Value = "0.0541"
$ws.Range("E42").Value = "  +11.13%  "

$ws.Range("D43").Value = "0.594"
$ws.Range("E43").Value = "  +9.33%  "

$ws.Range("D44").Value = "1.27"
$ws.Range("E44").Value = "  +13.68%  "

$ws.Range("B45").Value = "WhiteBITCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D45").Value = "10.13"
$ws.Range("E45").Value = "  +0.18%  "

$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").Value = "4.64"
$ws.Range("E46").Value = "  +23.96%  "

$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D47").Value = "254.76"
$ws.Range("E47").Value = "  +33.37%  "

$ws.Range("D48").Value = "0.0887"
$ws.Range("E48").Value = "  +12.54%  "

$ws.Range("D49").Value = "0.0220"
$ws.Range("E49").Value = "  +10.44%  "

$ws.Range("D50").Value = "1.896.78"
$ws.Range("E50").Value = "  +3.77%  "

$ws.Range("D51").Value = "16.91"
$ws.Range("E51").Value = "  +9.68%  "
